# Updates the 100 arithmetic-answer cells in the single table (20 rows x
# 5 columns) of the document, in row-major order, to match the "output
# generated at c8c62b6" commit. We assign Range.Text directly (after
# trimming the cell-end mark with MoveEnd) rather than using
# Range.Find.Execute, because this runtime's Find/Replace matches by
# text content rather than strict range position, which corrupts results
# when the same old value (e.g. "20-5=15") occurs in more than one cell.
$d = $word.ActiveDocument
$newValues = @(
    "17+9=26",
    "9+12=21",
    "50-14=36",
    "86+5=91",
    "34+38=72",
    "67+28=95",
    "72+9=81",
    "81-62=19",
    "37+57=94",
    "54+9=63",
    "91-7=84",
    "4+67=71",
    "80-15=65",
    "85-49=36",
    "9+8=17",
    "18+8=26",
    "25+39=64",
    "11-8=3",
    "70-48=22",
    "18+48=66",
    "66-9=57",
    "39+19=58",
    "36+8=44",
    "39+18=57",
    "45+8=53",
    "6+56=62",
    "45+47=92",
    "31-12=19",
    "82-15=67",
    "69+26=95",
    "79+17=96",
    "94-47=47",
    "91-48=43",
    "35+36=71",
    "73-38=35",
    "72-56=16",
    "92-17=75",
    "19+32=51",
    "38+7=45",
    "70-4=66",
    "18+58=76",
    "8+69=77",
    "19+39=58",
    "41-9=32",
    "58+18=76",
    "29+64=93",
    "64-16=48",
    "96-47=49",
    "67+27=94",
    "84-58=26",
    "7+9=16",
    "37+46=83",
    "91-25=66",
    "60-8=52",
    "6+5=11",
    "7+14=21",
    "76-47=29",
    "25+39=64",
    "69+25=94",
    "89+3=92",
    "81-68=13",
    "7+26=33",
    "65-36=29",
    "29+5=34",
    "62-24=38",
    "93-44=49",
    "19+46=65",
    "34+47=81",
    "7+26=33",
    "14+59=73",
    "55+27=82",
    "73-56=17",
    "85-19=66",
    "55+8=63",
    "72-59=13",
    "61-46=15",
    "72+9=81",
    "17+56=73",
    "91-35=56",
    "85-78=7",
    "29+29=58",
    "7+7=14",
    "3+38=41",
    "19+17=36",
    "8+49=57",
    "93-17=76",
    "21-12=9",
    "50-33=17",
    "7+6=13",
    "35+46=81",
    "90-24=66",
    "32-17=15",
    "6+37=43",
    "8+57=65",
    "32-19=13",
    "39+45=84",
    "69+17=86",
    "25+48=73",
    "21-7=14",
    "65+19=84",
)

$t = $d.Tables.Item(1)
$idx = 0
foreach ($row in $t.Rows) {
    foreach ($cell in $row.Cells) {
        if ($idx -lt $newValues.Count) {
            $r = $cell.Range
            $r.MoveEnd(1, -1) | Out-Null
            $r.Text = $newValues[$idx]
        }
        $idx = $idx + 1
    }
}
Write-Host "Replaced $idx cells"